# Update database values and change read_price algorithm
# (EPS based on latest capital, row 27, columns D:H)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("D27").Value = 100
$ws.Range("E27").Value = 136
$ws.Range("F27").Value = 196
$ws.Range("G27").Value = 502
$ws.Range("H27").Value = 930
